$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Pan's bake shop" row (currently row 5) down to row 6,
# and put the new "Amorino" entry into row 5.

$ws.Range("C6").Value = "Pan's bake shop"
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = 6
$ws.Range("H6").Value = 5

$ws.Range("C5").Value = "Amorino"
$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = 8
$ws.Range("H5").Value = 8

# Update the active selection to match the saved view state.
$ws.Range("E7").Select()
